$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell B1 becomes the text label "AB"
$ws.Range("B1").Value = "AB"

# Columns A (year) and B (value) for rows 2..452
$arr = New-Object "object[,]" 451,2
$arr[0,0] = 1600
$arr[0,1] = 18259.18871024664
$arr[1,0] = 1601
$arr[1,1] = 21.99771500693703
$arr[2,0] = 1602
$arr[2,1] = 22.35501739751943
$arr[3,0] = 1603
$arr[3,1] = 22.72018747415721
$arr[4,0] = 1604
$arr[4,1] = 23.09349612067184
$arr[5,0] = 1605
$arr[5,1] = 23.47522521919411
$arr[6,0] = 1606
$arr[6,1] = 23.86566797991999
$arr[7,0] = 1607
$arr[7,1] = 24.26512926851581
$arr[8,0] = 1608
$arr[8,1] = 24.67392593006804
$arr[9,0] = 1609
$arr[9,1] = 25.09238710808253
$arr[10,0] = 1610
$arr[10,1] = 25.52085455752959
$arr[11,0] = 1611
$arr[11,1] = 25.95968295059534
$arr[12,0] = 1612
$arr[12,1] = 26.40924017339148
$arr[13,0] = 1613
$arr[13,1] = 26.86990761259826
$arr[14,0] = 1614
$arr[14,1] = 27.34208043001852
$arr[15,0] = 1615
$arr[15,1] = 27.82616782384388
$arr[16,0] = 1616
$arr[16,1] = 28.32259327456793
$arr[17,0] = 1617
$arr[17,1] = 28.83179477403698
$arr[18,0] = 1618
$arr[18,1] = 29.35422503577518
$arr[19,0] = 1619
$arr[19,1] = 29.89035168453679
$arr[20,0] = 1620
$arr[20,1] = 30.4406574235003
$arr[21,0] = 1621
$arr[21,1] = 31.005640176754
$arr[22,0] = 1622
$arr[22,1] = 31.58581320528901
$arr[23,0] = 1623
$arr[23,1] = 32.18170519445266
$arr[24,0] = 1624
$arr[24,1] = 32.79386031053316
$arr[25,0] = 1625
$arr[25,1] = 33.42283822468092
$arr[26,0] = 1626
$arr[26,1] = 34.06921410174754
$arr[27,0] = 1627
$arr[27,1] = 34.73357855202004
$arr[28,0] = 1628
$arr[28,1] = 35.41653754375653
$arr[29,0] = 1629
$arr[29,1] = 36.11871227415051
$arr[30,0] = 1630
$arr[30,1] = 36.84073899679978
$arr[31,0] = 1631
$arr[31,1] = 37.58326880346247
$arr[32,0] = 1632
$arr[32,1] = 38.34696735797368
$arr[33,0] = 1633
$arr[33,1] = 39.13251458038019
$arr[34,0] = 1634
$arr[34,1] = 39.94060427916613
$arr[35,0] = 1635
$arr[35,1] = 40.77194372980077
$arr[36,0] = 1636
$arr[36,1] = 41.62725319751743
$arr[37,0] = 1637
$arr[37,1] = 42.50726540284681
$arr[38,0] = 1638
$arr[38,1] = 43.41272492801646
$arr[39,0] = 1639
$arr[39,1] = 44.34438756276081
$arr[40,0] = 1640
$arr[40,1] = 45.30301958797527
$arr[41,0] = 1641
$arr[41,1] = 46.28939699628943
$arr[42,0] = 1642
$arr[42,1] = 47.30430464786993
$arr[43,0] = 1643
$arr[43,1] = 48.34853536101957
$arr[44,0] = 1644
$arr[44,1] = 49.42288893642074
$arr[45,0] = 1645
$arr[45,1] = 50.52817111448525
$arr[46,0] = 1646
$arr[46,1] = 51.66519246568422
$arr[47,0] = 1647
$arr[47,1] = 52.83476721321507
$arr[48,0] = 1648
$arr[48,1] = 54.03771198839962
$arr[49,0] = 1649
$arr[49,1] = 55.27484451907639
$arr[50,0] = 1650
$arr[50,1] = 56.54698225131961
$arr[51,0] = 1651
$arr[51,1] = 57.85494090554295
$arr[52,0] = 1652
$arr[52,1] = 59.19953296782885
$arr[53,0] = 1653
$arr[53,1] = 60.58156611819201
$arr[54,0] = 1654
$arr[54,1] = 62.00184159698965
$arr[55,0] = 1655
$arr[55,1] = 63.46115251178984
$arr[56,0] = 1656
$arr[56,1] = 64.96028208697609
$arr[57,0] = 1657
$arr[57,1] = 66.50000185831631
$arr[58,0] = 1658
$arr[58,1] = 68.08106981580327
$arr[59,0] = 1659
$arr[59,1] = 69.70422849785771
$arr[60,0] = 1660
$arr[60,1] = 71.37020304036398
$arr[61,0] = 1661
$arr[61,1] = 73.07969918444853
$arr[62,0] = 1662
$arr[62,1] = 74.83340124748987
$arr[63,0] = 1663
$arr[63,1] = 76.631970061288
$arr[64,0] = 1664
$arr[64,1] = 78.47604088306278
$arr[65,0] = 1665
$arr[65,1] = 80.36622128364627
$arr[66,0] = 1666
$arr[66,1] = 82.3030890189283
$arr[67,0] = 1667
$arr[67,1] = 84.2871898901698
$arr[68,0] = 1668
$arr[68,1] = 86.31903559915621
$arr[69,0] = 1669
$arr[69,1] = 88.39910160460747
$arr[70,0] = 1670
$arr[70,1] = 90.52782498667966
$arr[71,0] = 1671
$arr[71,1] = 92.70560232603259
$arr[72,0] = 1672
$arr[72,1] = 94.93278760497715
$arr[73,0] = 1673
$arr[73,1] = 97.20969013778004
$arr[74,0] = 1674
$arr[74,1] = 99.53657253773689
$arr[75,0] = 1675
$arr[75,1] = 101.9136487286799
$arr[76,0] = 1676
$arr[76,1] = 104.3410820088526
$arr[77,0] = 1677
$arr[77,1] = 106.81898317493
$arr[78,0] = 1678
$arr[78,1] = 109.3474087144767
$arr[79,0] = 1679
$arr[79,1] = 111.9263590747855
$arr[80,0] = 1680
$arr[80,1] = 114.5557770163513
$arr[81,0] = 1681
$arr[81,1] = 117.2355460591294
$arr[82,0] = 1682
$arr[82,1] = 119.9654890297508
$arr[83,0] = 1683
$arr[83,1] = 122.7453667177182
$arr[84,0] = 1684
$arr[84,1] = 125.5748766486934
$arr[85,0] = 1685
$arr[85,1] = 128.453651982429
$arr[86,0] = 1686
$arr[86,1] = 131.3812605435118
$arr[87,0] = 1687
$arr[87,1] = 134.3572039921193
$arr[88,0] = 1688
$arr[88,1] = 137.3809171417924
$arr[89,0] = 1689
$arr[89,1] = 140.4517674318333
$arr[90,0] = 1690
$arr[90,1] = 143.5690545604813
$arr[91,0] = 1691
$arr[91,1] = 146.7320102852427
$arr[92,0] = 1692
$arr[92,1] = 149.9397983965323
$arr[93,0] = 1693
$arr[93,1] = 153.1915148699795
$arr[94,0] = 1694
$arr[94,1] = 156.4861882023152
$arr[95,0] = 1695
$arr[95,1] = 159.8227799356874
$arr[96,0] = 1696
$arr[96,1] = 163.2001853742121
$arr[97,0] = 1697
$arr[97,1] = 166.6172344963555
$arr[98,0] = 1698
$arr[98,1] = 170.0726930660104
$arr[99,0] = 1699
$arr[99,1] = 173.5652639444205
$arr[100,0] = 1700
$arr[100,1] = 177.0935886049974
$arr[101,0] = 1701
$arr[101,1] = 180.6562488515202
$arr[102,0] = 1702
$arr[102,1] = 184.2517687402285
$arr[103,0] = 1703
$arr[103,1] = 187.8786167057306
$arr[104,0] = 1704
$arr[104,1] = 191.5352078889212
$arr[105,0] = 1705
$arr[105,1] = 195.2199066658106
$arr[106,0] = 1706
$arr[106,1] = 198.9310293740265
$arr[107,0] = 1707
$arr[107,1] = 202.6668472341436
$arr[108,0] = 1708
$arr[108,1] = 206.4255894615041
$arr[109,0] = 1709
$arr[109,1] = 210.205446563647
$arr[110,0] = 1710
$arr[110,1] = 214.0045738179703
$arr[111,0] = 1711
$arr[111,1] = 217.8210949229457
$arr[112,0] = 1712
$arr[112,1] = 221.6531058161913
$arr[113,0] = 1713
$arr[113,1] = 225.4986786512572
$arr[114,0] = 1714
$arr[114,1] = 229.355865924767
$arr[115,0] = 1715
$arr[115,1] = 233.2227047446877
$arr[116,0] = 1716
$arr[116,1] = 237.0972212299049
$arr[117,0] = 1717
$arr[117,1] = 240.9774350303924
$arr[118,0] = 1718
$arr[118,1] = 244.861363957112
$arr[119,0] = 1719
$arr[119,1] = 248.7470287097589
$arr[120,0] = 1720
$arr[120,1] = 252.6324576902162
$arr[121,0] = 1721
$arr[121,1] = 256.5156918888818
$arr[122,0] = 1722
$arr[122,1] = 260.3947898307531
$arr[123,0] = 1723
$arr[123,1] = 264.2678325676109
$arr[124,0] = 1724
$arr[124,1] = 268.1329287023707
$arr[125,0] = 1725
$arr[125,1] = 271.9882194311769
$arr[126,0] = 1726
$arr[126,1] = 275.8318835889193
$arr[127,0] = 1727
$arr[127,1] = 279.6621426830574
$arr[128,0] = 1728
$arr[128,1] = 283.4772659013049
$arr[129,0] = 1729
$arr[129,1] = 287.2755750775008
$arr[130,0] = 1730
$arr[130,1] = 291.0554496013349
$arr[131,0] = 1731
$arr[131,1] = 294.8153312563945
$arr[132,0] = 1732
$arr[132,1] = 298.5537289715328
$arr[133,0] = 1733
$arr[133,1] = 302.2692234713658
$arr[134,0] = 1734
$arr[134,1] = 305.9604718102494
$arr[135,0] = 1735
$arr[135,1] = 309.6262117764566
$arr[136,0] = 1736
$arr[136,1] = 313.2652661518181
$arr[137,0] = 1737
$arr[137,1] = 316.8765468132468
$arr[138,0] = 1738
$arr[138,1] = 320.4590586633055
$arr[139,0] = 1739
$arr[139,1] = 324.0119033766492
$arr[140,0] = 1740
$arr[140,1] = 327.5342829502903
$arr[141,0] = 1741
$arr[141,1] = 331.0255030463709
$arr[142,0] = 1742
$arr[142,1] = 334.4849761159714
$arr[143,0] = 1743
$arr[143,1] = 337.9122242943128
$arr[144,0] = 1744
$arr[144,1] = 341.3068820572338
$arr[145,0] = 1745
$arr[145,1] = 344.668698630938
$arr[146,0] = 1746
$arr[146,1] = 347.9975401462906
$arr[147,0] = 1747
$arr[147,1] = 351.2933915313643
$arr[148,0] = 1748
$arr[148,1] = 354.5563581357482
$arr[149,0] = 1749
$arr[149,1] = 357.7866670811355
$arr[150,0] = 1750
$arr[150,1] = 360.9846683345498
$arr[151,0] = 1751
$arr[151,1] = 364.1508355001117
$arr[152,0] = 1752
$arr[152,1] = 367.2857663276872
$arr[153,0] = 1753
$arr[153,1] = 370.3901829364555
$arr[154,0] = 1754
$arr[154,1] = 373.4649317532704
$arr[155,0] = 1755
$arr[155,1] = 376.5109831667239
$arr[156,0] = 1756
$arr[156,1] = 379.5294308977003
$arr[157,0] = 1757
$arr[157,1] = 382.5214910900263
$arr[158,0] = 1758
$arr[158,1] = 385.4885011239331
$arr[159,0] = 1759
$arr[159,1] = 388.4319181573987
$arr[160,0] = 1760
$arr[160,1] = 391.3533174006109
$arr[161,0] = 1761
$arr[161,1] = 394.254390130016
$arr[162,0] = 1762
$arr[162,1] = 397.1369414493696
$arr[163,0] = 1763
$arr[163,1] = 400.0028878055717
$arr[164,0] = 1764
$arr[164,1] = 402.8542542688058
$arr[165,0] = 1765
$arr[165,1] = 405.6931715865001
$arr[166,0] = 1766
$arr[166,1] = 408.5218730215203
$arr[167,0] = 1767
$arr[167,1] = 411.3426909860804
$arr[168,0] = 1768
$arr[168,1] = 414.1580534833135
$arr[169,0] = 1769
$arr[169,1] = 416.9704803687838
$arr[170,0] = 1770
$arr[170,1] = 419.7825794452664
$arr[171,0] = 1771
$arr[171,1] = 422.5970424043135
$arr[172,0] = 1772
$arr[172,1] = 425.4166406287086
$arr[173,0] = 1773
$arr[173,1] = 428.2442208701277
$arr[174,0] = 1774
$arr[174,1] = 431.0827008170092
$arr[175,0] = 1775
$arr[175,1] = 433.9350645676213
$arr[176,0] = 1776
$arr[176,1] = 436.8043580235385
$arr[177,0] = 1777
$arr[177,1] = 439.6936842191691
$arr[178,0] = 1778
$arr[178,1] = 442.6061986026952
$arr[179,0] = 1779
$arr[179,1] = 445.5451042842436
$arr[180,0] = 1780
$arr[180,1] = 448.5136472668016
$arr[181,0] = 1781
$arr[181,1] = 451.5151116753077
$arr[182,0] = 1782
$arr[182,1] = 454.5528149993492
$arr[183,0] = 1783
$arr[183,1] = 457.6301033649266
$arr[184,0] = 1784
$arr[184,1] = 460.7503468498226
$arr[185,0] = 1785
$arr[185,1] = 463.9169348575126
$arr[186,0] = 1786
$arr[186,1] = 467.133271564106
$arr[187,0] = 1787
$arr[187,1] = 470.4027714516816
$arr[188,0] = 1788
$arr[188,1] = 473.7288549423801
$arr[189,0] = 1789
$arr[189,1] = 477.1149441454932
$arr[190,0] = 1790
$arr[190,1] = 480.5644587303741
$arr[191,0] = 1791
$arr[191,1] = 484.0808119375876
$arr[192,0] = 1792
$arr[192,1] = 487.6674067387199
$arr[193,0] = 1793
$arr[193,1] = 491.3276321563722
$arr[194,0] = 1794
$arr[194,1] = 495.0648597543096
$arr[195,0] = 1795
$arr[195,1] = 498.882440306941
$arr[196,0] = 1796
$arr[196,1] = 502.78370065707
$arr[197,0] = 1797
$arr[197,1] = 506.7719407700533
$arr[198,0] = 1798
$arr[198,1] = 510.8504309913395
$arr[199,0] = 1799
$arr[199,1] = 515.0224095144214
$arr[200,0] = 1800
$arr[200,1] = 519.2910800647168
$arr[201,0] = 1801
$arr[201,1] = 523.6597266565124
$arr[202,0] = 1802
$arr[202,1] = 528.1312572580513
$arr[203,0] = 1803
$arr[203,1] = 532.7088654344733
$arr[204,0] = 1804
$arr[204,1] = 537.3955984211483
$arr[205,0] = 1805
$arr[205,1] = 542.1944608000078
$arr[206,0] = 1806
$arr[206,1] = 547.1084135603541
$arr[207,0] = 1807
$arr[207,1] = 552.1403733959695
$arr[208,0] = 1808
$arr[208,1] = 557.2932122382902
$arr[209,0] = 1809
$arr[209,1] = 562.5697570250419
$arr[210,0] = 1810
$arr[210,1] = 567.9727897030446
$arr[211,0] = 1811
$arr[211,1] = 573.5050474629827
$arr[212,0] = 1812
$arr[212,1] = 579.1692232037598
$arr[213,0] = 1813
$arr[213,1] = 584.9679662234774
$arr[214,0] = 1814
$arr[214,1] = 590.903883132731
$arr[215,0] = 1815
$arr[215,1] = 596.9795389864368
$arr[216,0] = 1816
$arr[216,1] = 603.1974586294024
$arr[217,0] = 1817
$arr[217,1] = 609.5601282500593
$arr[218,0] = 1818
$arr[218,1] = 616.0699971369842
$arr[219,0] = 1819
$arr[219,1] = 622.7294796325876
$arr[220,0] = 1820
$arr[220,1] = 629.5409572757866
$arr[221,0] = 1821
$arr[221,1] = 636.5067811294068
$arr[222,0] = 1822
$arr[222,1] = 643.62927428263
$arr[223,0] = 1823
$arr[223,1] = 650.9107345233597
$arr[224,0] = 1824
$arr[224,1] = 658.3534371702063
$arr[225,0] = 1825
$arr[225,1] = 665.9596380594409
$arr[226,0] = 1826
$arr[226,1] = 673.7315766754909
$arr[227,0] = 1827
$arr[227,1] = 681.6714794190318
$arr[228,0] = 1828
$arr[228,1] = 689.7815630030281
$arr[229,0] = 1829
$arr[229,1] = 698.0640379690873
$arr[230,0] = 1830
$arr[230,1] = 706.5211123140078
$arr[231,0] = 1831
$arr[231,1] = 715.1549952193088
$arr[232,0] = 1832
$arr[232,1] = 723.9679008743735
$arr[233,0] = 1833
$arr[233,1] = 732.9620523835915
$arr[234,0] = 1834
$arr[234,1] = 742.1396857500092
$arr[235,0] = 1835
$arr[235,1] = 751.5030539254782
$arr[236,0] = 1836
$arr[236,1] = 761.0544309188067
$arr[237,0] = 1837
$arr[237,1] = 770.7961159532938
$arr[238,0] = 1838
$arr[238,1] = 780.7304376644787
$arr[239,0] = 1839
$arr[239,1] = 790.8597583289458
$arr[240,0] = 1840
$arr[240,1] = 801.1864781165597
$arr[241,0] = 1841
$arr[241,1] = 811.7130393561529
$arr[242,0] = 1842
$arr[242,1] = 822.4419308062517
$arr[243,0] = 1843
$arr[243,1] = 833.3756919230296
$arr[244,0] = 1844
$arr[244,1] = 844.516917115375
$arr[245,0] = 1845
$arr[245,1] = 855.8682599799839
$arr[246,0] = 1846
$arr[246,1] = 867.4324375066946
$arr[247,0] = 1847
$arr[247,1] = 879.212234246151
$arr[248,0] = 1848
$arr[248,1] = 891.2105064310077
$arr[249,0] = 1849
$arr[249,1] = 903.4301860424165
$arr[250,0] = 1850
$arr[250,1] = 915.8742848123463
$arr[251,0] = 1851
$arr[251,1] = 928.5458981544316
$arr[252,0] = 1852
$arr[252,1] = 941.4482090137712
$arr[253,0] = 1853
$arr[253,1] = 954.5844916266868
$arr[254,0] = 1854
$arr[254,1] = 967.9581151827665
$arr[255,0] = 1855
$arr[255,1] = 981.5725473794575
$arr[256,0] = 1856
$arr[256,1] = 995.4313578601016
$arr[257,0] = 1857
$arr[257,1] = 1009.538221526589
$arr[258,0] = 1858
$arr[258,1] = 1023.896921717743
$arr[259,0] = 1859
$arr[259,1] = 1038.5113532424
$arr[260,0] = 1860
$arr[260,1] = 1053.385525259791
$arr[261,0] = 1861
$arr[261,1] = 1068.523563995009
$arr[262,0] = 1862
$arr[262,1] = 1083.929715280489
$arr[263,0] = 1863
$arr[263,1] = 1099.608346914475
$arr[264,0] = 1864
$arr[264,1] = 1115.563950821763
$arr[265,0] = 1865
$arr[265,1] = 1131.801145010901
$arr[266,0] = 1866
$arr[266,1] = 1148.324675312737
$arr[267,0] = 1867
$arr[267,1] = 1165.13941688965
$arr[268,0] = 1868
$arr[268,1] = 1182.25037550478
$arr[269,0] = 1869
$arr[269,1] = 1199.662688536991
$arr[270,0] = 1870
$arr[270,1] = 1217.38162573076
$arr[271,0] = 1871
$arr[271,1] = 1235.412589666172
$arr[272,0] = 1872
$arr[272,1] = 1253.761115936995
$arr[273,0] = 1873
$arr[273,1] = 1272.432873022394
$arr[274,0] = 1874
$arr[274,1] = 1291.433661837643
$arr[275,0] = 1875
$arr[275,1] = 1310.769414950497
$arr[276,0] = 1876
$arr[276,1] = 1330.446195446806
$arr[277,0] = 1877
$arr[277,1] = 1350.470195431066
$arr[278,0] = 1878
$arr[278,1] = 1370.847734145452
$arr[279,0] = 1879
$arr[279,1] = 1391.585255692044
$arr[280,0] = 1880
$arr[280,1] = 1412.68932634109
$arr[281,0] = 1881
$arr[281,1] = 1434.166631409149
$arr[282,0] = 1882
$arr[282,1] = 1456.023971689826
$arr[283,0] = 1883
$arr[283,1] = 1478.268259420072
$arr[284,0] = 1884
$arr[284,1] = 1500.906513765048
$arr[285,0] = 1885
$arr[285,1] = 1523.94585580254
$arr[286,0] = 1886
$arr[286,1] = 1547.393502991415
$arr[287,0] = 1887
$arr[287,1] = 1571.25676310455
$arr[288,0] = 1888
$arr[288,1] = 1595.543027609674
$arr[289,0] = 1889
$arr[289,1] = 1620.259764480044
$arr[290,0] = 1890
$arr[290,1] = 1645.414510417427
$arr[291,0] = 1891
$arr[291,1] = 1671.014862470294
$arr[292,0] = 1892
$arr[292,1] = 1697.068469031354
$arr[293,0] = 1893
$arr[293,1] = 1723.583020195332
$arr[294,0] = 1894
$arr[294,1] = 1750.566237464467
$arr[295,0] = 1895
$arr[295,1] = 1778.025862784213
$arr[296,0] = 1896
$arr[296,1] = 1805.969646895187
$arr[297,0] = 1897
$arr[297,1] = 1834.40533698781
$arr[298,0] = 1898
$arr[298,1] = 1863.340663648054
$arr[299,0] = 1899
$arr[299,1] = 1892.783327081497
$arr[300,0] = 1900
$arr[300,1] = 1922.740982605587
$arr[301,0] = 1901
$arr[301,1] = 1953.221225404408
$arr[302,0] = 1902
$arr[302,1] = 1984.231574535401
$arr[303,0] = 1903
$arr[303,1] = 2015.7794561849
$arr[304,0] = 1904
$arr[304,1] = 2047.872186170569
$arr[305,0] = 1905
$arr[305,1] = 2080.516951688013
$arr[306,0] = 1906
$arr[306,1] = 2113.720792303831
$arr[307,0] = 1907
$arr[307,1] = 2147.490580201277
$arr[308,0] = 1908
$arr[308,1] = 2181.832999683855
$arr[309,0] = 1909
$arr[309,1] = 2216.75452594625
$arr[310,0] = 1910
$arr[310,1] = 2252.261403129742
$arr[311,0] = 1911
$arr[311,1] = 2288.359621676602
$arr[312,0] = 1912
$arr[312,1] = 2325.05489500325
$arr[313,0] = 1913
$arr[313,1] = 2362.352635520108
$arr[314,0] = 1914
$arr[314,1] = 2400.257930025001
$arr[315,0] = 1915
$arr[315,1] = 2438.775514502697
$arr[316,0] = 1916
$arr[316,1] = 2477.909748371586
$arr[317,0] = 1917
$arr[317,1] = 2517.664588215964
$arr[318,0] = 1918
$arr[318,1] = 2558.043561053834
$arr[319,0] = 1919
$arr[319,1] = 2599.049737194707
$arr[320,0] = 1920
$arr[320,1] = 2640.685702740131
$arr[321,0] = 1921
$arr[321,1] = 2682.953531796343
$arr[322,0] = 1922
$arr[322,1] = 2725.854758461723
$arr[323,0] = 1923
$arr[323,1] = 2769.390348670602
$arr[324,0] = 1924
$arr[324,1] = 2813.560671966328
$arr[325,0] = 1925
$arr[325,1] = 2858.365473294011
$arr[326,0] = 1926
$arr[326,1] = 2903.80384490365
$arr[327,0] = 1927
$arr[327,1] = 2949.87419845893
$arr[328,0] = 1928
$arr[328,1] = 2996.574237456018
$arr[329,0] = 1929
$arr[329,1] = 3043.900930060042
$arr[330,0] = 1930
$arr[330,1] = 3091.850482472806
$arr[331,0] = 1931
$arr[331,1] = 3140.418312948511
$arr[332,0] = 1932
$arr[332,1] = 3189.599026583601
$arr[333,0] = 1933
$arr[333,1] = 3239.38639100589
$arr[334,0] = 1934
$arr[334,1] = 3289.773313096231
$arr[335,0] = 1935
$arr[335,1] = 3340.751816879131
$arr[336,0] = 1936
$arr[336,1] = 3392.313022719929
$arr[337,0] = 1937
$arr[337,1] = 3444.447127970283
$arr[338,0] = 1938
$arr[338,1] = 3497.143389210858
$arr[339,0] = 1939
$arr[339,1] = 3550.390106231271
$arr[340,0] = 1940
$arr[340,1] = 3604.174607898279
$arr[341,0] = 1941
$arr[341,1] = 3658.483240060958
$arr[342,0] = 1942
$arr[342,1] = 3713.301355637705
$arr[343,0] = 1943
$arr[343,1] = 3768.613307036059
$arr[344,0] = 1944
$arr[344,1] = 3824.402441043935
$arr[345,0] = 1945
$arr[345,1] = 3880.6510963448
$arr[346,0] = 1946
$arr[346,1] = 3937.340603784142
$arr[347,0] = 1947
$arr[347,1] = 3994.451289534381
$arr[348,0] = 1948
$arr[348,1] = 4051.962481275933
$arr[349,0] = 1949
$arr[349,1] = 4109.852517529426
$arr[350,0] = 1950
$arr[350,1] = 4168.098760246334
$arr[351,0] = 1951
$arr[351,1] = 4226.6776107731
$arr[352,0] = 1952
$arr[352,1] = 4285.56452928343
$arr[353,0] = 1953
$arr[353,1] = 4344.734057772562
$arr[354,0] = 1954
$arr[354,1] = 4404.159846689085
$arr[355,0] = 1955
$arr[355,1] = 4463.81468527281
$arr[356,0] = 1956
$arr[356,1] = 5137.183946017288
$arr[357,0] = 1957
$arr[357,1] = 5205.353158556394
$arr[358,0] = 1958
$arr[358,1] = 5273.68432244705
$arr[359,0] = 1959
$arr[359,1] = 5342.143055014871
$arr[360,0] = 1960
$arr[360,1] = 5410.694232568432
$arr[361,0] = 1961
$arr[361,1] = 5479.302050825856
$arr[362,0] = 1962
$arr[362,1] = 5547.930089636425
$arr[363,0] = 1963
$arr[363,1] = 5616.541381922724
$arr[364,0] = 1964
$arr[364,1] = 5685.0984867664
$arr[365,0] = 1965
$arr[365,1] = 5753.563566539692
$arr[366,0] = 1966
$arr[366,1] = 5821.898467947105
$arr[367,0] = 1967
$arr[367,1] = 5890.064806848937
$arr[368,0] = 1968
$arr[368,1] = 5958.024056688463
$arr[369,0] = 1969
$arr[369,1] = 6025.737640348148
$arr[370,0] = 1970
$arr[370,1] = 6093.167025229895
$arr[371,0] = 1971
$arr[371,1] = 5210.736942727913
$arr[372,0] = 1972
$arr[372,1] = 5267.194849410677
$arr[373,0] = 1973
$arr[373,1] = 5323.315650654282
$arr[374,0] = 1974
$arr[374,1] = 5379.067755683577
$arr[375,0] = 1975
$arr[375,1] = 5434.420047258678
$arr[376,0] = 1976
$arr[376,1] = 5489.341973977868
$arr[377,0] = 1977
$arr[377,1] = 5543.803643687975
$arr[378,0] = 1978
$arr[378,1] = 5597.775917722701
$arr[379,0] = 1979
$arr[379,1] = 5651.230505681922
$arr[380,0] = 1980
$arr[380,1] = 5704.140060454534
$arr[381,0] = 1981
$arr[381,1] = 4203.838755945564
$arr[382,0] = 1982
$arr[382,1] = 4241.624661000586
$arr[383,0] = 1983
$arr[383,1] = 4278.957451482866
$arr[384,0] = 1984
$arr[384,1] = 4315.820601560365
$arr[385,0] = 1985
$arr[385,1] = 4352.198616112737
$arr[386,0] = 1986
$arr[386,1] = 4388.077096488729
$arr[387,0] = 1987
$arr[387,1] = 4423.442804863361
$arr[388,0] = 1988
$arr[388,1] = 4458.283726972249
$arr[389,0] = 1989
$arr[389,1] = 4492.589133000035
$arr[390,0] = 1990
$arr[390,1] = 4526.349636409002
$arr[391,0] = 1991
$arr[391,1] = 7400.161573106768
$arr[392,0] = 1992
$arr[392,1] = 7453.149589756364
$arr[393,0] = 1993
$arr[393,1] = 7505.221504614012
$arr[394,0] = 1994
$arr[394,1] = 7556.371679846349
$arr[395,0] = 1995
$arr[395,1] = 7606.597073188725
$arr[396,0] = 1996
$arr[396,1] = 7655.897306411346
$arr[397,0] = 1997
$arr[397,1] = 7704.274728221711
$arr[398,0] = 1998
$arr[398,1] = 7751.734471310482
$arr[399,0] = 1999
$arr[399,1] = 7798.284503291008
$arr[400,0] = 2000
$arr[400,1] = 7843.935671249782
$arr[401,0] = 2001
$arr[401,1] = 12089.7757232225
$arr[402,0] = 2002
$arr[402,1] = 12157.05081430843
$arr[403,0] = 2003
$arr[403,1] = 12223.02523408918
$arr[404,0] = 2004
$arr[404,1] = 12287.73373191133
$arr[405,0] = 2005
$arr[405,1] = 12351.21565166777
$arr[406,0] = 2006
$arr[406,1] = 12413.51493363417
$arr[407,0] = 2007
$arr[407,1] = 12474.68010401605
$arr[408,0] = 2008
$arr[408,1] = 12534.76425180913
$arr[409,0] = 2009
$arr[409,1] = 12593.82499267761
$arr[410,0] = 2010
$arr[410,1] = 12651.92441944137
$arr[411,0] = 2011
$arr[411,1] = 12105.70066544271
$arr[412,0] = 2012
$arr[412,1] = 12159.4043721004
$arr[413,0] = 2013
$arr[413,1] = 12212.39475946571
$arr[414,0] = 2014
$arr[414,1] = 12264.74729213655
$arr[415,0] = 2015
$arr[415,1] = 12316.54126018691
$arr[416,0] = 2016
$arr[416,1] = 12367.85964806489
$arr[417,0] = 2017
$arr[417,1] = 12418.78898815572
$arr[418,0] = 2018
$arr[418,1] = 12469.4191985934
$arr[419,0] = 2019
$arr[419,1] = 12519.84340491348
$arr[420,0] = 2020
$arr[420,1] = 12570.15774511293
$arr[421,0] = 2021
$arr[421,1] = 12620.46115772751
$arr[422,0] = 2022
$arr[422,1] = 12670.85515248061
$arr[423,0] = 2023
$arr[423,1] = 12721.44356314544
$arr[424,0] = 2024
$arr[424,1] = 12772.33228219719
$arr[425,0] = 2025
$arr[425,1] = 12823.62897693279
$arr[426,0] = 2026
$arr[426,1] = 12875.44278673031
$arr[427,0] = 2027
$arr[427,1] = 12927.88400116044
$arr[428,0] = 2028
$arr[428,1] = 12981.06371875367
$arr[429,0] = 2029
$arr[429,1] = 13035.09348626017
$arr[430,0] = 2030
$arr[430,1] = 13090.08491833099
$arr[431,0] = 2031
$arr[431,1] = 13146.14929766479
$arr[432,0] = 2032
$arr[432,1] = 13203.39715573274
$arr[433,0] = 2033
$arr[433,1] = 13261.93783437344
$arr[434,0] = 2034
$arr[434,1] = 13321.87902862448
$arr[435,0] = 2035
$arr[435,1] = 13383.32631139385
$arr[436,0] = 2036
$arr[436,1] = 13446.38264066379
$arr[437,0] = 2037
$arr[437,1] = 13511.14785019487
$arr[438,0] = 2038
$arr[438,1] = 13577.71812482126
$arr[439,0] = 2039
$arr[439,1] = 13646.1854617259
$arr[440,0] = 2040
$arr[440,1] = 13716.6371192694
$arr[441,0] = 2041
$arr[441,1] = 13789.15505520801
$arr[442,0] = 2042
$arr[442,1] = 13863.81535642161
$arr[443,0] = 2043
$arr[443,1] = 13940.68766251234
$arr[444,0] = 2044
$arr[444,1] = 14019.83458593257
$arr[445,0] = 2045
$arr[445,1] = 14101.31113157673
$arr[446,0] = 2046
$arr[446,1] = 14185.1641190611
$arr[447,0] = 2047
$arr[447,1] = 14271.43161122209
$arr[448,0] = 2048
$arr[448,1] = 14360.14235257009
$arr[449,0] = 2049
$arr[449,1] = 14451.31522181803
$arr[450,0] = 2050
$arr[450,1] = 14544.95870275859

$ws.Range("A2:B452").Value = $arr
